# PRISMA.pptx — "removed papers that do not fit our criteria"
#
# 1) datetimeFigureOut field on the slide master + all 11 slide layouts:
#    "12/18/20" -> "1/11/21"
# 2) Flow-diagram box text: "(N = 177)" -> "(N = 179)"
# 3) Flow-diagram box text: "(N = 37)"  -> "(N = 35)" (ends up as two runs:
#    "(N " / "= 35)")

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" date field wherever it appears: the
#    slide master and every custom (slide) layout each carry their own
#    "Date Placeholder" shape holding the field. (No recursion is used here
#    — this runtime's PowerShell engine blows its expression-nesting budget
#    on recursive functions, so everything below is a flat loop.)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "12/18/20") {
        $sh.TextFrame.TextRange.Text = "1/11/21"
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "12/18/20") {
            $sh.TextFrame.TextRange.Text = "1/11/21"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) & 3) The two counts live in text boxes nested one level inside group
#    shapes on slide 1 of the PRISMA flow diagram.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$box177 = $null
$box37 = $null

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $top = $slide.Shapes.Item($i)
    if ($top.Type -eq 6) {
        for ($j = 1; $j -le $top.GroupItems.Count; $j++) {
            $inner = $top.GroupItems.Item($j)
            if ($inner.Name -eq "TextBox 38") {
                $box177 = $inner
            }
            if ($inner.Name -eq "TextBox 42") {
                $box37 = $inner
            }
        }
    }
}

# "Incorrect or missing experimental design elements (N = 177)" -> (N = 179)
$tr177 = $box177.TextFrame.TextRange
$full177 = $tr177.Text
$idx177 = $full177.IndexOf("(N = 177)")
$sub177 = $tr177.Characters($idx177 + 1, 9)
$sub177.Text = "(N = 179)"

# "Papers included in meta-analysis (N = 37)" -> split into "(N " + "= 35)"
$tr37 = $box37.TextFrame.TextRange
$full37 = $tr37.Text
$idx37 = $full37.IndexOf("(N = 37)")
# Only touch the digits + closing paren ("37)" -> "35)") so the engine
# re-flows the tail into its own run, matching "(N " / "= 35)".
$subNum = $tr37.Characters($idx37 + 4, 5)
$subNum.Text = "= 35)"
